# Scheduled-runner refresh of cached crafting/leve profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per updated market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1515.0416  # H12
$ws.Cells.Item(12, 9).Value = 969.8095  # I12
$ws.Cells.Item(12, 10).Value = 5331.6665  # J12
$ws.Cells.Item(12, 11).Value = 969.8095  # K12
$ws.Cells.Item(12, 12).Value = 5331.6665  # L12
$ws.Cells.Item(12, 13).Value = -799.8095  # M12
$ws.Cells.Item(12, 14).Value = -5671.6665  # N12
$ws.Cells.Item(70, 8).Value = 2015.5385  # H70
$ws.Cells.Item(70, 9).Value = 1962.2858  # I70
$ws.Cells.Item(70, 10).Value = 2077.6667  # J70
$ws.Cells.Item(70, 11).Value = 5886.857400000001  # K70
$ws.Cells.Item(70, 12).Value = 6233.000100000001  # L70
$ws.Cells.Item(70, 13).Value = -5616.857400000001  # M70
$ws.Cells.Item(70, 14).Value = -6773.000100000001  # N70
$ws.Cells.Item(73, 8).Value = 2015.5385  # H73
$ws.Cells.Item(73, 9).Value = 1962.2858  # I73
$ws.Cells.Item(73, 10).Value = 2077.6667  # J73
$ws.Cells.Item(73, 11).Value = 5886.857400000001  # K73
$ws.Cells.Item(73, 12).Value = 6233.000100000001  # L73
$ws.Cells.Item(73, 13).Value = -4950.857400000001  # M73
$ws.Cells.Item(73, 14).Value = -8105.000100000001  # N73
$ws.Cells.Item(100, 8).Value = 3889.4443  # H100
$ws.Cells.Item(100, 9).Value = 3975.625  # I100
$ws.Cells.Item(100, 10).Value = 3200  # J100
$ws.Cells.Item(100, 11).Value = 3975.625  # K100
$ws.Cells.Item(100, 12).Value = 3200  # L100
$ws.Cells.Item(100, 13).Value = -3434.625  # M100
$ws.Cells.Item(100, 14).Value = -4282  # N100
$ws.Cells.Item(107, 8).Value = 348.66666  # H107
$ws.Cells.Item(107, 9).Value = 348.66666  # I107
$ws.Cells.Item(107, 10).Value = 0  # J107
$ws.Cells.Item(107, 11).Value = 348.66666  # K107
$ws.Cells.Item(107, 12).Value = 0  # L107
$ws.Cells.Item(107, 13).ClearContents()  # M107
$ws.Cells.Item(107, 14).Value = 1571.33334  # N107
$ws.Cells.Item(116, 8).Value = 2664.3333  # H116
$ws.Cells.Item(116, 9).Value = 1997  # I116
$ws.Cells.Item(116, 11).Value = 1997  # K116
$ws.Cells.Item(116, 13).Value = 1445  # M116

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(50, 8).Value = 44000  # H50
$ws.Cells.Item(50, 10).Value = 44000  # J50
$ws.Cells.Item(50, 12).Value = 44000  # L50
$ws.Cells.Item(50, 14).Value = -45428  # N50
$ws.Cells.Item(61, 8).Value = 7290.125  # H61
$ws.Cells.Item(61, 9).Value = 7617.2856  # I61
$ws.Cells.Item(61, 10).Value = 5000  # J61
$ws.Cells.Item(61, 11).Value = 7617.2856  # K61
$ws.Cells.Item(61, 12).Value = 5000  # L61
$ws.Cells.Item(61, 13).Value = -7405.2856  # M61
$ws.Cells.Item(61, 14).Value = -5424  # N61
$ws.Cells.Item(74, 8).Value = 2014.25  # H74
$ws.Cells.Item(74, 9).Value = 852.5  # I74
$ws.Cells.Item(74, 11).Value = 852.5  # K74
$ws.Cells.Item(74, 13).Value = 21.5  # M74
$ws.Cells.Item(77, 8).Value = 2014.25  # H77
$ws.Cells.Item(77, 9).Value = 852.5  # I77
$ws.Cells.Item(77, 11).Value = 4262.5  # K77
$ws.Cells.Item(77, 13).Value = 105.5  # M77
$ws.Cells.Item(108, 8).Value = 58562.855  # H108
$ws.Cells.Item(108, 10).Value = 58562.855  # J108
$ws.Cells.Item(108, 12).Value = 58562.855  # L108
$ws.Cells.Item(108, 14).Value = -66242.85500000001  # N108
$ws.Cells.Item(132, 8).Value = 1300.0769  # H132
$ws.Cells.Item(132, 9).Value = 1175.0834  # I132
$ws.Cells.Item(132, 11).Value = 3525.2502  # K132
$ws.Cells.Item(132, 13).Value = -995.2501999999999  # M132
$ws.Cells.Item(136, 8).Value = 7290.125  # H136
$ws.Cells.Item(136, 9).Value = 7617.2856  # I136
$ws.Cells.Item(136, 10).Value = 5000  # J136
$ws.Cells.Item(136, 11).Value = 22851.8568  # K136
$ws.Cells.Item(136, 12).Value = 15000  # L136
$ws.Cells.Item(136, 13).Value = -20301.8568  # M136
$ws.Cells.Item(136, 14).Value = -20100  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(32, 8).Value = 20029  # H32
$ws.Cells.Item(32, 10).Value = 20029  # J32
$ws.Cells.Item(32, 12).Value = 20029  # L32
$ws.Cells.Item(32, 14).Value = -20797  # N32
$ws.Cells.Item(37, 8).Value = 921.5  # H37
$ws.Cells.Item(37, 9).Value = 1080  # I37
$ws.Cells.Item(37, 10).Value = 129  # J37
$ws.Cells.Item(37, 11).Value = 1080  # K37
$ws.Cells.Item(37, 12).Value = 129  # L37
$ws.Cells.Item(37, 13).Value = -943  # M37
$ws.Cells.Item(37, 14).Value = -403  # N37
$ws.Cells.Item(99, 8).Value = 3150  # H99
$ws.Cells.Item(99, 9).Value = 6300  # I99
$ws.Cells.Item(99, 11).Value = 6300  # K99
$ws.Cells.Item(99, 13).Value = -4802  # M99
$ws.Cells.Item(132, 8).Value = 50129.168  # H132
$ws.Cells.Item(132, 10).Value = 50129.168  # J132
$ws.Cells.Item(132, 12).Value = 50129.168  # L132
$ws.Cells.Item(132, 14).Value = -60249.168  # N132

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(32, 8).Value = 809.8  # H32
$ws.Cells.Item(32, 9).Value = 809.8  # I32
$ws.Cells.Item(32, 11).Value = 809.8  # K32
$ws.Cells.Item(32, 13).Value = -493.8  # M32
$ws.Cells.Item(58, 8).Value = 3558.5715  # H58
$ws.Cells.Item(58, 9).Value = 3558.5715  # I58
$ws.Cells.Item(58, 10).Value = 0  # J58
$ws.Cells.Item(58, 11).Value = 3558.5715  # K58
$ws.Cells.Item(58, 12).Value = 0  # L58
$ws.Cells.Item(58, 13).ClearContents()  # M58
$ws.Cells.Item(58, 14).Value = -3355.5715  # N58
$ws.Cells.Item(134, 8).Value = 2336.5334  # H134
$ws.Cells.Item(134, 9).Value = 2378.9167  # I134
$ws.Cells.Item(134, 10).Value = 2167  # J134
$ws.Cells.Item(134, 11).Value = 7136.750100000001  # K134
$ws.Cells.Item(134, 12).Value = 6501  # L134
$ws.Cells.Item(134, 13).Value = -4601.750100000001  # M134
$ws.Cells.Item(134, 14).Value = -11571  # N134
$ws.Cells.Item(136, 8).Value = 3558.5715  # H136
$ws.Cells.Item(136, 9).Value = 3558.5715  # I136
$ws.Cells.Item(136, 10).Value = 0  # J136
$ws.Cells.Item(136, 11).Value = 10675.7145  # K136
$ws.Cells.Item(136, 12).Value = 0  # L136
$ws.Cells.Item(136, 13).ClearContents()  # M136
$ws.Cells.Item(136, 14).Value = -8125.7145  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 203.08333  # H38
$ws.Cells.Item(38, 9).Value = 50.4  # I38
$ws.Cells.Item(38, 10).Value = 312.14285  # J38
$ws.Cells.Item(38, 11).Value = 151.2  # K38
$ws.Cells.Item(38, 12).Value = 936.4285500000001  # L38
$ws.Cells.Item(38, 13).Value = 195.8  # M38
$ws.Cells.Item(38, 14).Value = -1630.42855  # N38
$ws.Cells.Item(103, 8).Value = 4500  # H103
$ws.Cells.Item(103, 9).Value = 0  # I103
$ws.Cells.Item(103, 10).Value = 4500  # J103
$ws.Cells.Item(103, 11).Value = 0  # K103
$ws.Cells.Item(103, 12).ClearContents()  # L103
$ws.Cells.Item(103, 13).Value = 13500  # M103
$ws.Cells.Item(103, 14).Value = -15258  # N103
$ws.Cells.Item(104, 8).Value = 15000  # H104
$ws.Cells.Item(104, 10).Value = 15000  # J104
$ws.Cells.Item(104, 12).Value = 45000  # L104
$ws.Cells.Item(104, 14).Value = -50242  # N104
$ws.Cells.Item(109, 8).Value = 508  # H109
$ws.Cells.Item(109, 9).Value = 405.5  # I109
$ws.Cells.Item(109, 10).Value = 713  # J109
$ws.Cells.Item(109, 11).Value = 1216.5  # K109
$ws.Cells.Item(109, 12).Value = 2139  # L109
$ws.Cells.Item(109, 13).Value = -176.5  # M109
$ws.Cells.Item(109, 14).Value = -4219  # N109
$ws.Cells.Item(128, 8).Value = 0  # H128
$ws.Cells.Item(128, 9).Value = 0  # I128
$ws.Cells.Item(128, 11).Value = 0  # K128
$ws.Cells.Item(128, 13).ClearContents()  # M128
$ws.Cells.Item(141, 8).Value = 11291.667  # H141
$ws.Cells.Item(141, 10).Value = 17877  # J141
$ws.Cells.Item(141, 12).Value = 53631  # L141
$ws.Cells.Item(141, 14).Value = -63991  # N141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 2000  # H107
$ws.Cells.Item(107, 9).Value = 0  # I107
$ws.Cells.Item(107, 10).Value = 2000  # J107
$ws.Cells.Item(107, 11).Value = 0  # K107
$ws.Cells.Item(107, 12).ClearContents()  # L107
$ws.Cells.Item(107, 13).Value = 2000  # M107
$ws.Cells.Item(107, 14).Value = -5840  # N107
$ws.Cells.Item(108, 8).Value = 61995.625  # H108
$ws.Cells.Item(108, 10).Value = 61995.625  # J108
$ws.Cells.Item(108, 12).Value = 61995.625  # L108
$ws.Cells.Item(108, 14).Value = -69675.625  # N108

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 999  # H68
$ws.Cells.Item(68, 9).Value = 999  # I68
$ws.Cells.Item(68, 11).Value = 999  # K68
$ws.Cells.Item(68, 13).Value = -250  # M68
$ws.Cells.Item(71, 8).Value = 999  # H71
$ws.Cells.Item(71, 9).Value = 999  # I71
$ws.Cells.Item(71, 11).Value = 4995  # K71
$ws.Cells.Item(71, 13).Value = -1251  # M71
$ws.Cells.Item(132, 8).Value = 10000  # H132
$ws.Cells.Item(132, 9).Value = 10000  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 30000  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 13).ClearContents()  # M132
$ws.Cells.Item(132, 14).Value = -27470  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(31, 8).Value = 19994  # H31
$ws.Cells.Item(31, 10).Value = 19994  # J31
$ws.Cells.Item(31, 12).Value = 19994  # L31
$ws.Cells.Item(31, 14).Value = -20690  # N31
$ws.Cells.Item(107, 8).Value = 448.8  # H107
$ws.Cells.Item(107, 9).Value = 453.16666  # I107
$ws.Cells.Item(107, 10).Value = 442.25  # J107
$ws.Cells.Item(107, 11).Value = 1359.49998  # K107
$ws.Cells.Item(107, 12).Value = 1326.75  # L107
$ws.Cells.Item(107, 13).Value = 560.5000199999999  # M107
$ws.Cells.Item(107, 14).Value = -5166.75  # N107
$ws.Cells.Item(126, 8).Value = 1416.3334  # H126
$ws.Cells.Item(126, 9).Value = 1374.5  # I126
$ws.Cells.Item(126, 10).Value = 1500  # J126
$ws.Cells.Item(126, 11).Value = 4123.5  # K126
$ws.Cells.Item(126, 12).Value = 4500  # L126
$ws.Cells.Item(126, 13).Value = -1653.5  # M126
$ws.Cells.Item(126, 14).Value = -9440  # N126
$ws.Cells.Item(132, 8).Value = 250537.75  # H132
$ws.Cells.Item(132, 9).Value = 250537.75  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 751613.25  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 13).ClearContents()  # M132
$ws.Cells.Item(132, 14).Value = -749083.25  # N132
